$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPS Data")

# Two swapped-out data values (commit #3473)
$ws.Range("L3").Value = 227440.2
$ws.Range("B6").Value = 22482006
$ws.Range("B10").Value = 22482007

# Update the sheet's current selection
$ws.Activate()
$ws.Range("L2:L10").Select()
